$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A21").Value = "小美人魚雨衣"
$ws.Range("B21").Value = 4212
$ws.Range("C21").Value = 2948
$ws.Range("D21").Formula = "=C21*D1"
$ws.Range("E21").Formula = "=D21+50"

$ws.Range("A22").Value = "奇奇蒂蒂好收納雨衣"
$ws.Range("B22").Value = 2700
$ws.Range("C22").Value = 2160
$ws.Range("D22").Formula = "=C22*D1"
$ws.Range("E22").Formula = "=D22+50"

$ws.Range("E22").Select()
